$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 03:52"

# Update per-country statistics (columns B..H) for the affected rows.
# Each row below keeps the same country name in column A; only the
# numeric columns (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) change.

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

Set-Row 18  @(12232, 49, 127, 11539, 296, 2, 566)
Set-Row 20  @(10331, 47, 6694, 3445, 55, 6, 192)
Set-Row 42  @(2439, 296, 633, 1681, 89, 31, 125)
Set-Row 60  @(1160, 54, 241, 918, 4, 0, 1)
Set-Row 61  @(1120, 0, 81, 959, 1, 0, 80)
Set-Row 62  @(1108, 0, 62, 1027, 14, 0, 19)
Set-Row 99  @(287, 73, 31, 251, 2, 0, 5)
Set-Row 100 @(277, 0, 35, 210, 14, 0, 32)
Set-Row 101 @(254, 0, 24, 229, 0, 0, 1)
Set-Row 102 @(253, 0, 26, 217, 0, 0, 10)
Set-Row 103 @(245, 0, 95, 150, 8, 0, 0)
Set-Row 104 @(244, 0, 7, 230, 2, 0, 7)
Set-Row 105 @(241, 0, 5, 236, 3, 0, 0)
Set-Row 106 @(238, 0, 35, 198, 2, 0, 5)
Set-Row 107 @(233, 0, 1, 230, 4, 0, 2)
Set-Row 108 @(226, 0, 92, 132, 1, 0, 2)
Set-Row 109 @(216, 0, 33, 179, 5, 0, 4)
Set-Row 110 @(194, 11, 2, 178, 3, 3, 14)
Set-Row 111 @(188, 0, 39, 147, 6, 0, 2)
Set-Row 155 @(31, 0, 8, 18, 8, 1, 5)
